# "added new window for the items"
# Appends a new order row (row 7) to the "Order Data" sheet:
#   Order ID=10, Order Type="delivery", Items="[2, 2, 2, 2, 2, 2]",
#   Order Completed Status=FALSE, Order status="InProgress", Customer ID=4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "delivery"
$ws.Range("C7").Value = "[2, 2, 2, 2, 2, 2]"
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = "InProgress"
$ws.Range("F7").Value = 4
